$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two more levels done: append the new timing row (row 46) under the
# existing "D11:D45" shared-formula block.
$ws.Range("A46").Value = "Get Key"
$ws.Range("B46").Value = 79223
$ws.Range("C46").Value = 77419
$ws.Range("D46").Formula = "=C46-B46"

# Scroll the view down and move the active selection past the new row,
# matching where the author left off editing.
$win = $excel.ActiveWindow
$win.ScrollRow = 38
$win.ScrollColumn = 1
[void]$ws.Range("C47").Select()
